# Auto-generated COM-interop edit script.
# Adds 5 learners, 10 enrollments, 4 courses plus rebuilt aggregate sheets
# ("EnrollmentsWithCourses", "EnrollmentsByDate", "EnrollmentsByCourse",
# "CoursesByCategory", "EnrollmentsByCourseWithProgress") to match the new data.

$wb = $excel.ActiveWorkbook

# ===== Learners: append L006-L010 =====
$ws = $wb.Worksheets.Item("Learners")
$dateTextFmt = $ws.Cells.Item(1,1).NumberFormat  # "General" -> used as a base; reset below per-cell to "@"
$ws.Cells.Item(7,1).Value = "L006"
$ws.Cells.Item(7,2).Value = "frank@example.com"
$ws.Cells.Item(7,3).Value = "AU"
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "2024-03-20"
$ws.Cells.Item(8,1).Value = "L007"
$ws.Cells.Item(8,2).Value = "grace@example.com"
$ws.Cells.Item(8,3).Value = "US"
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "2024-03-22"
$ws.Cells.Item(9,1).Value = "L008"
$ws.Cells.Item(9,2).Value = "harry@example.com"
$ws.Cells.Item(9,3).Value = "UK"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "2024-03-25"
$ws.Cells.Item(10,1).Value = "L009"
$ws.Cells.Item(10,2).Value = "irene@example.com"
$ws.Cells.Item(10,3).Value = "IN"
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "2024-03-28"
$ws.Cells.Item(11,1).Value = "L010"
$ws.Cells.Item(11,2).Value = "jack@example.com"
$ws.Cells.Item(11,3).Value = "CA"
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "2024-04-01"

# ===== Enrollments: append E007-E016 (enroll_date stored as a real date serial) =====
$ws = $wb.Worksheets.Item("Enrollments")
$enrollDateFmt = $ws.Cells.Item(2,5).NumberFormat
$ws.Cells.Item(8,1).Value = "E007"
$ws.Cells.Item(8,2).Value = "L006"
$ws.Cells.Item(8,3).Value = "C101"
$ws.Cells.Item(8,4).Value = "I001"
$ws.Cells.Item(8,5).NumberFormat = $enrollDateFmt
$ws.Cells.Item(8,5).Value = 45369
$ws.Cells.Item(8,6).Value = 20
$ws.Cells.Item(8,7).Value = 60
$ws.Cells.Item(8,8).Value = $false
$ws.Cells.Item(9,1).Value = "E008"
$ws.Cells.Item(9,2).Value = "L006"
$ws.Cells.Item(9,3).Value = "C105"
$ws.Cells.Item(9,4).Value = "I003"
$ws.Cells.Item(9,5).NumberFormat = $enrollDateFmt
$ws.Cells.Item(9,5).Value = 45373
$ws.Cells.Item(9,6).Value = 40
$ws.Cells.Item(9,7).Value = 120
$ws.Cells.Item(9,8).Value = $false
$ws.Cells.Item(10,1).Value = "E009"
$ws.Cells.Item(10,2).Value = "L007"
$ws.Cells.Item(10,3).Value = "C102"
$ws.Cells.Item(10,4).Value = "I002"
$ws.Cells.Item(10,5).NumberFormat = $enrollDateFmt
$ws.Cells.Item(10,5).Value = 45374
$ws.Cells.Item(10,6).Value = 60
$ws.Cells.Item(10,7).Value = 180
$ws.Cells.Item(10,8).Value = $false
$ws.Cells.Item(11,1).Value = "E010"
$ws.Cells.Item(11,2).Value = "L007"
$ws.Cells.Item(11,3).Value = "C105"
$ws.Cells.Item(11,4).Value = "I003"
$ws.Cells.Item(11,5).NumberFormat = $enrollDateFmt
$ws.Cells.Item(11,5).Value = 45376
$ws.Cells.Item(11,6).Value = 10
$ws.Cells.Item(11,7).Value = 45
$ws.Cells.Item(11,8).Value = $false
$ws.Cells.Item(12,1).Value = "E011"
$ws.Cells.Item(12,2).Value = "L008"
$ws.Cells.Item(12,3).Value = "C103"
$ws.Cells.Item(12,4).Value = "I001"
$ws.Cells.Item(12,5).NumberFormat = $enrollDateFmt
$ws.Cells.Item(12,5).Value = 45377
$ws.Cells.Item(12,6).Value = 80
$ws.Cells.Item(12,7).Value = 300
$ws.Cells.Item(12,8).Value = $true
$ws.Cells.Item(13,1).Value = "E012"
$ws.Cells.Item(13,2).Value = "L009"
$ws.Cells.Item(13,3).Value = "C106"
$ws.Cells.Item(13,4).Value = "I004"
$ws.Cells.Item(13,5).NumberFormat = $enrollDateFmt
$ws.Cells.Item(13,5).Value = 45380
$ws.Cells.Item(13,6).Value = 30
$ws.Cells.Item(13,7).Value = 90
$ws.Cells.Item(13,8).Value = $false
$ws.Cells.Item(14,1).Value = "E013"
$ws.Cells.Item(14,2).Value = "L009"
$ws.Cells.Item(14,3).Value = "C107"
$ws.Cells.Item(14,4).Value = "I004"
$ws.Cells.Item(14,5).NumberFormat = $enrollDateFmt
$ws.Cells.Item(14,5).Value = 45381
$ws.Cells.Item(14,6).Value = 0
$ws.Cells.Item(14,7).Value = 0
$ws.Cells.Item(14,8).Value = $false
$ws.Cells.Item(15,1).Value = "E014"
$ws.Cells.Item(15,2).Value = "L010"
$ws.Cells.Item(15,3).Value = "C108"
$ws.Cells.Item(15,4).Value = "I004"
$ws.Cells.Item(15,5).NumberFormat = $enrollDateFmt
$ws.Cells.Item(15,5).Value = 45384
$ws.Cells.Item(15,6).Value = 15
$ws.Cells.Item(15,7).Value = 75
$ws.Cells.Item(15,8).Value = $false
$ws.Cells.Item(16,1).Value = "E015"
$ws.Cells.Item(16,2).Value = "L010"
$ws.Cells.Item(16,3).Value = "C101"
$ws.Cells.Item(16,4).Value = "I001"
$ws.Cells.Item(16,5).NumberFormat = $enrollDateFmt
$ws.Cells.Item(16,5).Value = 45387
$ws.Cells.Item(16,6).Value = 50
$ws.Cells.Item(16,7).Value = 200
$ws.Cells.Item(16,8).Value = $false
$ws.Cells.Item(17,1).Value = "E016"
$ws.Cells.Item(17,2).Value = "L005"
$ws.Cells.Item(17,3).Value = "C105"
$ws.Cells.Item(17,4).Value = "I003"
$ws.Cells.Item(17,5).NumberFormat = $enrollDateFmt
$ws.Cells.Item(17,5).Value = 45369
$ws.Cells.Item(17,6).Value = 30
$ws.Cells.Item(17,7).Value = 100
$ws.Cells.Item(17,8).Value = $false

# ===== Courses: append C105-C108 =====
$ws = $wb.Worksheets.Item("Courses")
$ws.Cells.Item(6,1).Value = "C105"
$ws.Cells.Item(6,2).Value = "Data Visualization with Power BI"
$ws.Cells.Item(6,3).Value = "Data Science"
$ws.Cells.Item(6,4).Value = "Intermediate"
$ws.Cells.Item(6,5).Value = 420
$ws.Cells.Item(7,1).Value = "C106"
$ws.Cells.Item(7,2).Value = "Advanced SQL Optimization"
$ws.Cells.Item(7,3).Value = "Data Science"
$ws.Cells.Item(7,4).Value = "Advanced"
$ws.Cells.Item(7,5).Value = 480
$ws.Cells.Item(8,1).Value = "C107"
$ws.Cells.Item(8,2).Value = "Cloud Data Warehousing with Snowflake"
$ws.Cells.Item(8,3).Value = "Data Engineering"
$ws.Cells.Item(8,4).Value = "Intermediate"
$ws.Cells.Item(8,5).Value = 540
$ws.Cells.Item(9,1).Value = "C108"
$ws.Cells.Item(9,2).Value = "Real-time Analytics with Kafka & Spark"
$ws.Cells.Item(9,3).Value = "Data Engineering"
$ws.Cells.Item(9,4).Value = "Advanced"
$ws.Cells.Item(9,5).Value = 600

# ===== EnrollmentsWithCourses: append E007-E016 joined with course info (enroll_date stays text) =====
$ws = $wb.Worksheets.Item("EnrollmentsWithCourses")
$ws.Cells.Item(8,1).Value = "E007"
$ws.Cells.Item(8,2).Value = "L006"
$ws.Cells.Item(8,3).Value = "C101"
$ws.Cells.Item(8,4).Value = "I001"
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = "2024-03-18"
$ws.Cells.Item(8,6).Value = 20
$ws.Cells.Item(8,7).Value = 60
$ws.Cells.Item(8,8).Value = $false
$ws.Cells.Item(8,9).Value = "Data Engineering Fundamentals"
$ws.Cells.Item(8,10).Value = "Data Science"
$ws.Cells.Item(8,11).Value = "Beginner"
$ws.Cells.Item(8,12).Value = 480
$ws.Cells.Item(9,1).Value = "E008"
$ws.Cells.Item(9,2).Value = "L006"
$ws.Cells.Item(9,3).Value = "C105"
$ws.Cells.Item(9,4).Value = "I003"
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value = "2024-03-22"
$ws.Cells.Item(9,6).Value = 40
$ws.Cells.Item(9,7).Value = 120
$ws.Cells.Item(9,8).Value = $false
$ws.Cells.Item(9,9).Value = "Data Visualization with Power BI"
$ws.Cells.Item(9,10).Value = "Data Science"
$ws.Cells.Item(9,11).Value = "Intermediate"
$ws.Cells.Item(9,12).Value = 420
$ws.Cells.Item(10,1).Value = "E009"
$ws.Cells.Item(10,2).Value = "L007"
$ws.Cells.Item(10,3).Value = "C102"
$ws.Cells.Item(10,4).Value = "I002"
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value = "2024-03-23"
$ws.Cells.Item(10,6).Value = 60
$ws.Cells.Item(10,7).Value = 180
$ws.Cells.Item(10,8).Value = $false
$ws.Cells.Item(10,9).Value = "Machine Learning Basics"
$ws.Cells.Item(10,10).Value = "Data Science"
$ws.Cells.Item(10,11).Value = "Intermediate"
$ws.Cells.Item(10,12).Value = 600
$ws.Cells.Item(11,1).Value = "E010"
$ws.Cells.Item(11,2).Value = "L007"
$ws.Cells.Item(11,3).Value = "C105"
$ws.Cells.Item(11,4).Value = "I003"
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value = "2024-03-25"
$ws.Cells.Item(11,6).Value = 10
$ws.Cells.Item(11,7).Value = 45
$ws.Cells.Item(11,8).Value = $false
$ws.Cells.Item(11,9).Value = "Data Visualization with Power BI"
$ws.Cells.Item(11,10).Value = "Data Science"
$ws.Cells.Item(11,11).Value = "Intermediate"
$ws.Cells.Item(11,12).Value = 420
$ws.Cells.Item(12,1).Value = "E011"
$ws.Cells.Item(12,2).Value = "L008"
$ws.Cells.Item(12,3).Value = "C103"
$ws.Cells.Item(12,4).Value = "I001"
$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,5).Value = "2024-03-26"
$ws.Cells.Item(12,6).Value = 80
$ws.Cells.Item(12,7).Value = 300
$ws.Cells.Item(12,8).Value = $true
$ws.Cells.Item(12,9).Value = "SQL for Analytics"
$ws.Cells.Item(12,10).Value = "Data Science"
$ws.Cells.Item(12,11).Value = "Beginner"
$ws.Cells.Item(12,12).Value = 360
$ws.Cells.Item(13,1).Value = "E012"
$ws.Cells.Item(13,2).Value = "L009"
$ws.Cells.Item(13,3).Value = "C106"
$ws.Cells.Item(13,4).Value = "I004"
$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value = "2024-03-29"
$ws.Cells.Item(13,6).Value = 30
$ws.Cells.Item(13,7).Value = 90
$ws.Cells.Item(13,8).Value = $false
$ws.Cells.Item(13,9).Value = "Advanced SQL Optimization"
$ws.Cells.Item(13,10).Value = "Data Science"
$ws.Cells.Item(13,11).Value = "Advanced"
$ws.Cells.Item(13,12).Value = 480
$ws.Cells.Item(14,1).Value = "E013"
$ws.Cells.Item(14,2).Value = "L009"
$ws.Cells.Item(14,3).Value = "C107"
$ws.Cells.Item(14,4).Value = "I004"
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value = "2024-03-30"
$ws.Cells.Item(14,6).Value = 0
$ws.Cells.Item(14,7).Value = 0
$ws.Cells.Item(14,8).Value = $false
$ws.Cells.Item(14,9).Value = "Cloud Data Warehousing with Snowflake"
$ws.Cells.Item(14,10).Value = "Data Engineering"
$ws.Cells.Item(14,11).Value = "Intermediate"
$ws.Cells.Item(14,12).Value = 540
$ws.Cells.Item(15,1).Value = "E014"
$ws.Cells.Item(15,2).Value = "L010"
$ws.Cells.Item(15,3).Value = "C108"
$ws.Cells.Item(15,4).Value = "I004"
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value = "2024-04-02"
$ws.Cells.Item(15,6).Value = 15
$ws.Cells.Item(15,7).Value = 75
$ws.Cells.Item(15,8).Value = $false
$ws.Cells.Item(15,9).Value = "Real-time Analytics with Kafka & Spark"
$ws.Cells.Item(15,10).Value = "Data Engineering"
$ws.Cells.Item(15,11).Value = "Advanced"
$ws.Cells.Item(15,12).Value = 600
$ws.Cells.Item(16,1).Value = "E015"
$ws.Cells.Item(16,2).Value = "L010"
$ws.Cells.Item(16,3).Value = "C101"
$ws.Cells.Item(16,4).Value = "I001"
$ws.Cells.Item(16,5).NumberFormat = "@"
$ws.Cells.Item(16,5).Value = "2024-04-05"
$ws.Cells.Item(16,6).Value = 50
$ws.Cells.Item(16,7).Value = 200
$ws.Cells.Item(16,8).Value = $false
$ws.Cells.Item(16,9).Value = "Data Engineering Fundamentals"
$ws.Cells.Item(16,10).Value = "Data Science"
$ws.Cells.Item(16,11).Value = "Beginner"
$ws.Cells.Item(16,12).Value = 480
$ws.Cells.Item(17,1).Value = "E016"
$ws.Cells.Item(17,2).Value = "L005"
$ws.Cells.Item(17,3).Value = "C105"
$ws.Cells.Item(17,4).Value = "I003"
$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,5).Value = "2024-03-18"
$ws.Cells.Item(17,6).Value = 30
$ws.Cells.Item(17,7).Value = 100
$ws.Cells.Item(17,8).Value = $false
$ws.Cells.Item(17,9).Value = "Data Visualization with Power BI"
$ws.Cells.Item(17,10).Value = "Data Science"
$ws.Cells.Item(17,11).Value = "Intermediate"
$ws.Cells.Item(17,12).Value = 420

# ===== EnrollmentsByDate: append new date/count rows =====
$ws = $wb.Worksheets.Item("EnrollmentsByDate")
$byDateFmt = $ws.Cells.Item(2,1).NumberFormat
$ws.Cells.Item(8,1).NumberFormat = $byDateFmt
$ws.Cells.Item(8,1).Value = 45369
$ws.Cells.Item(8,2).Value = 2
$ws.Cells.Item(9,1).NumberFormat = $byDateFmt
$ws.Cells.Item(9,1).Value = 45373
$ws.Cells.Item(9,2).Value = 1
$ws.Cells.Item(10,1).NumberFormat = $byDateFmt
$ws.Cells.Item(10,1).Value = 45374
$ws.Cells.Item(10,2).Value = 1
$ws.Cells.Item(11,1).NumberFormat = $byDateFmt
$ws.Cells.Item(11,1).Value = 45376
$ws.Cells.Item(11,2).Value = 1
$ws.Cells.Item(12,1).NumberFormat = $byDateFmt
$ws.Cells.Item(12,1).Value = 45377
$ws.Cells.Item(12,2).Value = 1
$ws.Cells.Item(13,1).NumberFormat = $byDateFmt
$ws.Cells.Item(13,1).Value = 45380
$ws.Cells.Item(13,2).Value = 1
$ws.Cells.Item(14,1).NumberFormat = $byDateFmt
$ws.Cells.Item(14,1).Value = 45381
$ws.Cells.Item(14,2).Value = 1
$ws.Cells.Item(15,1).NumberFormat = $byDateFmt
$ws.Cells.Item(15,1).Value = 45384
$ws.Cells.Item(15,2).Value = 1
$ws.Cells.Item(16,1).NumberFormat = $byDateFmt
$ws.Cells.Item(16,1).Value = 45387
$ws.Cells.Item(16,2).Value = 1

# ===== EnrollmentsByCourse: replace with the recomputed, alphabetised rollup (4 -> 7 rows) =====
$ws = $wb.Worksheets.Item("EnrollmentsByCourse")
$ws.Cells.Item(2,1).Value = "Advanced SQL Optimization"
$ws.Cells.Item(2,2).Value = 1
$ws.Cells.Item(3,1).Value = "Cloud Data Warehousing with Snowflake"
$ws.Cells.Item(3,2).Value = 1
$ws.Cells.Item(4,1).Value = "Data Engineering Fundamentals"
$ws.Cells.Item(4,2).Value = 5
$ws.Cells.Item(5,1).Value = "Data Visualization with Power BI"
$ws.Cells.Item(5,2).Value = 3
$ws.Cells.Item(6,1).Value = "Machine Learning Basics"
$ws.Cells.Item(6,2).Value = 3
$ws.Cells.Item(7,1).Value = "Real-time Analytics with Kafka & Spark"
$ws.Cells.Item(7,2).Value = 1
$ws.Cells.Item(8,1).Value = "SQL for Analytics"
$ws.Cells.Item(8,2).Value = 2

# ===== CoursesByCategory: counts change, rows stay the same (2 rows) =====
$ws = $wb.Worksheets.Item("CoursesByCategory")
$ws.Cells.Item(2,2).Value = 3
$ws.Cells.Item(3,2).Value = 5

# ===== EnrollmentsByCourseWithProgress: replace with recomputed, alphabetised rollup (4 -> 7 rows) =====
$ws = $wb.Worksheets.Item("EnrollmentsByCourseWithProgress")
$ws.Cells.Item(2,1).Value = "Advanced SQL Optimization"
$ws.Cells.Item(2,2).Value = 1
$ws.Cells.Item(2,3).Value = 30
$ws.Cells.Item(3,1).Value = "Cloud Data Warehousing with Snowflake"
$ws.Cells.Item(3,2).Value = 1
$ws.Cells.Item(3,3).Value = 0
$ws.Cells.Item(4,1).Value = "Data Engineering Fundamentals"
$ws.Cells.Item(4,2).Value = 5
$ws.Cells.Item(4,3).Value = 64
$ws.Cells.Item(5,1).Value = "Data Visualization with Power BI"
$ws.Cells.Item(5,2).Value = 3
$ws.Cells.Item(5,3).Value = 26.66666666666667
$ws.Cells.Item(6,1).Value = "Machine Learning Basics"
$ws.Cells.Item(6,2).Value = 3
$ws.Cells.Item(6,3).Value = 45
$ws.Cells.Item(7,1).Value = "Real-time Analytics with Kafka & Spark"
$ws.Cells.Item(7,2).Value = 1
$ws.Cells.Item(7,3).Value = 15
$ws.Cells.Item(8,1).Value = "SQL for Analytics"
$ws.Cells.Item(8,2).Value = 2
$ws.Cells.Item(8,3).Value = 52.5

